# Clasificador.xlsx - reorder indicators alphabetically within each
# category (column B) block, reclassify SFE-ICC from "ConfCons" into
# "Construc", and update the view selection, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clasificador")

# SFE-ICC (row 4) moves out of the standalone "ConfCons" category and
# into "Construc" so the whole 8-row block (rows 4:11) is re-sorted
# together alphabetically by the indicator code in column A.
$ws.Cells.Item(4, 2).Value2 = "Construc"

# Re-sort each category block by column A (ascending), mirroring the
# per-category alphabetical ordering applied in the workbook.
$ws.Range("A4:B11").Sort($ws.Range("A4:A11"), 1)
$ws.Range("A12:B22").Sort($ws.Range("A12:A22"), 1)
$ws.Range("A23:B28").Sort($ws.Range("A23:A28"), 1)
$ws.Range("A29:B36").Sort($ws.Range("A29:A36"), 1)
$ws.Range("A37:B42").Sort($ws.Range("A37:A42"), 1)
$ws.Range("A44:B64").Sort($ws.Range("A44:A64"), 1)
$ws.Range("A65:B68").Sort($ws.Range("A65:A68"), 1)
$ws.Range("A69:B76").Sort($ws.Range("A69:A76"), 1)
$ws.Range("A77:B80").Sort($ws.Range("A77:A80"), 1)
$ws.Range("A81:B84").Sort($ws.Range("A81:A84"), 1)

# Update the saved selection to match the new view position.
$ws.Range("B69").Select()
